# Updating sorting order check for utilitysummary
# - adds a new "prodfix" worksheet (PROD_Fix scenarios 1-4)
# - flips which worksheet tab / cell is active
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # NewImportLogic
$ws2 = $wb.Worksheets.Item(2)   # OldImportLogic

# ---------------------------------------------------------------------------
# 1. Add the new "prodfix" worksheet after the existing sheets
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "prodfix"

# ---------------------------------------------------------------------------
# 2. Populate the header row + 4 PROD_Fix scenario blocks.
#    Cells are written in the same order the strings were first introduced
#    so that newly-interned shared strings line up with the source edit.
# ---------------------------------------------------------------------------
$ws3.Range("A1").Value = "Name"
$ws3.Range("B1").Value = "Import_Pop"
$ws3.Range("H1").Value = "ExtractionFile"

$ws3.Range("A2").Value = "scenario1"
$ws3.Range("G2").Value = "\Testdata\Templates\UtilityOutcome\QOL_ECON_Staging\PROD_Fix\UtilityOutcome_Feature_Extraction_file_ExpectedData.xlsx"
$ws3.Range("H2").Value = "\Testdata\Templates\UtilityOutcome\QOL_ECON_Staging\PROD_Fix\Templates\UtilityOutcome_Feature_Extraction_file_QoL_UtilityData_ECON_NoUtility.xlsx"

$ws3.Range("A5").Value = "scenario2"
$ws3.Range("H5").Value = "\Testdata\Templates\UtilityOutcome\QOL_ECON_Staging\PROD_Fix\Templates\UtilityOutcome_Feature_Extraction_file_ECON_UtilityData_QoL_NoUtility.xlsx"

$ws3.Range("A8").Value = "scenario3"
$ws3.Range("H8").Value = "\Testdata\Templates\UtilityOutcome\QOL_ECON_Staging\PROD_Fix\Templates\UtilityOutcome_Feature_Extraction_file_Both_QoL_ECON_Utility.xlsx"

$ws3.Range("A11").Value = "scenario4"
$ws3.Range("H11").Value = "\Testdata\Templates\UtilityOutcome\QOL_ECON_Staging\PROD_Fix\Templates\UtilityOutcome_Feature_Extraction_file_NegativeScenario_QoL__ECON_NoUtility.xlsx"

$ws3.Range("C2").Value = "PRODFix_QOL_ECON - UtilityOutcome"
$ws3.Range("D2").Value = "PRODFix_QOL_ECON - UtilityOutcome_radio_button"
$ws3.Range("B2").Value = "UtilityOutcome - PRODFix_QOL_ECON - Ovid search - 9/19/2022"

$ws3.Range("K2").Value = "ExcelReport-PRODFix_QOL_ECON - UtilityOutcome-Quality of Life-"
$ws3.Range("K3").Value = "WordReport-PRODFix_QOL_ECON - UtilityOutcome-Quality of Life-"

# Remaining cells reuse strings already present in the shared string table.
$ws3.Range("C1").Value = "Population"
$ws3.Range("D1").Value = "Population_Radio_button"
$ws3.Range("E1").Value = "slrtype"
$ws3.Range("F1").Value = "slrtype_Radio_button"
$ws3.Range("G1").Value = "ExpectedSourceTemplateFile"
$ws3.Range("I1").Value = "ReportedVariables"
$ws3.Range("J1").Value = "Reportedvariable_checkbox"
$ws3.Range("K1").Value = "ExpectedFilenames"

$ws3.Range("E2").Value = "Quality of Life"
$ws3.Range("F2").Value = "Quality of Life_radio_button"
$ws3.Range("I2").Value = "reported_variable_section1"
$ws3.Range("J2").Value = "reported_variable_section1_checkbox"

$ws3.Range("A3").Value = "scenario1"
$ws3.Range("I3").Value = "reported_variable_section2"
$ws3.Range("J3").Value = "reported_variable_section2_checkbox"
$ws3.Range("K3").Value = "WordReport-PRODFix_QOL_ECON - UtilityOutcome-Quality of Life-"

$ws3.Range("K4").Value = "Report-"

$ws3.Range("B5").Value = "UtilityOutcome - PRODFix_QOL_ECON - Ovid search - 9/19/2022"
$ws3.Range("C5").Value = "PRODFix_QOL_ECON - UtilityOutcome"
$ws3.Range("D5").Value = "PRODFix_QOL_ECON - UtilityOutcome_radio_button"
$ws3.Range("E5").Value = "Quality of Life"
$ws3.Range("F5").Value = "Quality of Life_radio_button"
$ws3.Range("G5").Value = "\Testdata\Templates\UtilityOutcome\QOL_ECON_Staging\PROD_Fix\UtilityOutcome_Feature_Extraction_file_ExpectedData.xlsx"
$ws3.Range("I5").Value = "reported_variable_section1"
$ws3.Range("J5").Value = "reported_variable_section1_checkbox"

$ws3.Range("A6").Value = "scenario2"
$ws3.Range("I6").Value = "reported_variable_section2"
$ws3.Range("J6").Value = "reported_variable_section2_checkbox"

$ws3.Range("B8").Value = "UtilityOutcome - PRODFix_QOL_ECON - Ovid search - 9/19/2022"
$ws3.Range("C8").Value = "PRODFix_QOL_ECON - UtilityOutcome"
$ws3.Range("D8").Value = "PRODFix_QOL_ECON - UtilityOutcome_radio_button"
$ws3.Range("E8").Value = "Quality of Life"
$ws3.Range("F8").Value = "Quality of Life_radio_button"
$ws3.Range("G8").Value = "\Testdata\Templates\UtilityOutcome\QOL_ECON_Staging\PROD_Fix\UtilityOutcome_Feature_Extraction_file_ExpectedData.xlsx"
$ws3.Range("I8").Value = "reported_variable_section1"
$ws3.Range("J8").Value = "reported_variable_section1_checkbox"

$ws3.Range("A9").Value = "scenario3"
$ws3.Range("I9").Value = "reported_variable_section2"
$ws3.Range("J9").Value = "reported_variable_section2_checkbox"

$ws3.Range("B11").Value = "UtilityOutcome - PRODFix_QOL_ECON - Ovid search - 9/19/2022"
$ws3.Range("C11").Value = "PRODFix_QOL_ECON - UtilityOutcome"
$ws3.Range("D11").Value = "PRODFix_QOL_ECON - UtilityOutcome_radio_button"
$ws3.Range("E11").Value = "Quality of Life"
$ws3.Range("F11").Value = "Quality of Life_radio_button"
$ws3.Range("G11").Value = "\Testdata\Templates\UtilityOutcome\QOL_ECON_Staging\PROD_Fix\UtilityOutcome_Feature_Extraction_file_ExpectedData.xlsx"
$ws3.Range("I11").Value = "reported_variable_section1"
$ws3.Range("J11").Value = "reported_variable_section1_checkbox"

$ws3.Range("A12").Value = "scenario4"
$ws3.Range("I12").Value = "reported_variable_section2"
$ws3.Range("J12").Value = "reported_variable_section2_checkbox"

# ---------------------------------------------------------------------------
# 3. Column widths on the new sheet (matching the authored layout)
# ---------------------------------------------------------------------------
$ws3.Columns.Item(2).ColumnWidth = 27.72
$ws3.Columns.Item(3).ColumnWidth = 27.94
$ws3.Columns.Item(4).ColumnWidth = 39.94
$ws3.Columns.Item(5).ColumnWidth = 11.61
$ws3.Columns.Item(6).ColumnWidth = 23.61
$ws3.Columns.Item(7).ColumnWidth = 32.83
$ws3.Columns.Item(8).ColumnWidth = 26.94
$ws3.Columns.Item(9).ColumnWidth = 22.94
$ws3.Columns.Item(10).ColumnWidth = 31.83
$ws3.Columns.Item(11).ColumnWidth = 56.17

# ---------------------------------------------------------------------------
# 4. Print setup on the new sheet (portrait orientation was stamped when the
#    sheet was authored).
# ---------------------------------------------------------------------------
$ws3.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 5. Worksheet selection / active state.
#    Final state: "prodfix" keeps its own idle selection at A2; OldImportLogic
#    loses its tab-selected flag; NewImportLogic becomes the active tab with
#    C15 selected.
# ---------------------------------------------------------------------------
[void]$ws3.Range("A2").Select()

[void]$ws2.Activate()
[void]$ws2.Range("E2").Select()

[void]$ws1.Activate()
[void]$ws1.Range("C15").Select()
